# Refresh the cryptos price table (Coin/Link/Price/Volume) with the
# latest values from the scheduled GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.240.39"
$ws.Range("E2").Value = "  +2.92%  "

# Row 3
$ws.Range("D3").Value = "1.717.93"
$ws.Range("E3").Value = "  +3.33%  "

# Row 4
$ws.Range("D4").Value = "'0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "'239.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.00%  "

# Row 6
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.01%  "

# Row 7
$ws.Range("D7").Value = "'0.4730"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.49%  "

# Row 8
$ws.Range("D8").Value = "'0.2624"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").Value = "'0.06200"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.59%  "

# Row 10
$ws.Range("D10").Value = "1.713.29"
$ws.Range("E10").Value = "  +3.08%  "

# Row 11
$ws.Range("D11").Value = "'0.07071"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.51%  "

# Row 12
$ws.Range("D12").Value = "'15.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.58%  "

# Row 13
$ws.Range("D13").Value = "'0.5930"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.79%  "

# Row 14
$ws.Range("D14").Value = "'4.418"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.26%  "

# Row 15
$ws.Range("D15").Value = "'76.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.10%  "

# Row 16
$ws.Range("D16").Value = "'1.000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.01%  "

# Row 17
$ws.Range("D17").Value = "'1.000"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.04%  "

# Row 18
$ws.Range("D18").Value = "26.245.79"

# Row 19
$ws.Range("E19").Value = "  +0.60%  "

# Row 20
$ws.Range("D20").Value = "'11.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.78%  "

# Row 21
$ws.Range("D21").Value = "1.935.38"
$ws.Range("E21").Value = "  +3.48%  "

# Row 22
$ws.Range("D22").Value = "'4.555"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.83%  "

# Row 23
$ws.Range("D23").Value = "'8.717"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.59%  "

# Row 24
$ws.Range("D24").Value = "'5.267"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.03%  "

# Row 25
$ws.Range("D25").Value = "'134.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.31%  "

# Row 26
$ws.Range("D26").Value = "'15.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.00%  "

# Row 27
$ws.Range("D27").Value = "'1.401"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.17%  "

# Row 28
$ws.Range("D28").Value = "'1.762"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.29%  "

# Row 29
$ws.Range("D29").Value = "'107.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.19%  "

# Row 30
$ws.Range("D30").Value = "'3.967"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.88%  "

# Row 31
$ws.Range("D31").Value = "'3.677"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.63%  "

# Row 32
$ws.Range("D32").Value = "'0.07758"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.47%  "

# Row 33
$ws.Range("D33").Value = "'0.04431"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.90%  "

# Row 34
$ws.Range("D34").Value = "'2.615"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.06%  "

# Row 35
$ws.Range("D35").Value = "'0.9741"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.46%  "

# Row 36
$ws.Range("D36").Value = "'0.6168"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.18%  "

# Row 37
$ws.Range("D37").Value = "'0.9237"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.84%  "

# Row 38
$ws.Range("D38").Value = "'111.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +15.42%  "

# Row 39
$ws.Range("D39").Value = "'2.413"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.13%  "

# Row 40
$ws.Range("D40").Value = "'1.915"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.21%  "

# Row 41
$ws.Range("E41").Value = "  +0.06%  "

# Row 42
$ws.Range("D42").Value = "'0.01474"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.38%  "

# Row 43
$ws.Range("D43").Value = "'5.391"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +13.68%  "

# Row 44
$ws.Range("D44").Value = "'0.3810"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.51%  "

# Row 45
$ws.Range("D45").Value = "'0.1174"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.92%  "

# Row 46
$ws.Range("D46").Value = "'6.260"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.84%  "

# Row 47
$ws.Range("D47").Value = "'0.05258"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.13%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'7.754"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.08%  "

# Row 49
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'30.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.96%  "

# Row 50
$ws.Range("D50").Value = "'0.3377"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.75%  "

# Row 51
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.214"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.75%  "
